# Commit: Sat, May 16, 2020  2:04:59 AM
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in table style
#    {9EC5FFD9-EC63-46DD-A455-E93D7511953F}.
# 2) Swap the presentation's applied theme palette from the custom
#    "Integral" / "Red Violet" colours to the standard "Office Theme"
#    colours (the theme actually bound to the slide master).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{9EC5FFD9-EC63-46DD-A455-E93D7511953F}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colour swap --------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> "Office Theme" values
$officeColors = 0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
